# Insert a new "ID" column before the existing data (A..E -> B..F),
# and populate it with the sample identifiers for each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:E one place to the right, inserting a new column A.
$ws.Columns.Item(1).Insert()

# The new column A should look like the (now shifted) header column B,
# i.e. bold / centered / bordered header style.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A1").Value = "ID"

$ids = @(
    "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
    "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
    "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
